# excel tareas y BBDD
# Fill in assignee names + dates for the "Tareas" sheet task list, and
# flesh out the "home" task description.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tareas")

$xlPasteFormats = -4122

function Copy-Style {
    param($fromAddr, $toAddr)
    $ws.Range($fromAddr).Copy() | Out-Null
    $ws.Range($toAddr).PasteSpecial($xlPasteFormats) | Out-Null
}

# --- Update the "home" backend task description (A4) ---
$ws.Range("A4").Value = "home, crear ruta, cargarla, maquetar"

# --- Backend side: assignee (column G) + date (column H), matching style of F/A column ---
Copy-Style "F2" "G2"
$ws.Range("G2").Value = "alejandro"
Copy-Style "F3" "G3"
$ws.Range("G3").Value = "jonatan"
Copy-Style "F4" "G4"
$ws.Range("G4").Value = "pati"
Copy-Style "F5" "G5"
$ws.Range("G5").Value = "Edward"

# --- Frontend side: assignee (column B) + date (column C) ---
Copy-Style "A4" "B4"
$ws.Range("B4").Value = "Angel"
Copy-Style "A5" "B5"
$ws.Range("B5").Value = "fabrizo"

# --- Dates: copy style from an already-text-styled cell, then set value + number format ---
$dateCells = @(
    @{ Addr = "H2"; Style = "F2"; Value = 45861 },
    @{ Addr = "H3"; Style = "F3"; Value = 45862 },
    @{ Addr = "C4"; Style = "A4"; Value = 45861 },
    @{ Addr = "H4"; Style = "F4"; Value = 45863 },
    @{ Addr = "C5"; Style = "A5"; Value = 45861 },
    @{ Addr = "H5"; Style = "F5"; Value = 45864 },
    @{ Addr = "H6"; Style = "F6"; Value = 45865 },
    @{ Addr = "H7"; Style = "F7"; Value = 45866 },
    @{ Addr = "H8"; Style = "F8"; Value = 45867 }
)

foreach ($dc in $dateCells) {
    Copy-Style $dc.Style $dc.Addr
    $ws.Range($dc.Addr).Value = $dc.Value
    $ws.Range($dc.Addr).NumberFormat = "dd/mm/yyyy"
}

Write-Output "done"
